$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Mapping")

$ws.Range("A2").Value = -124.3795
$ws.Range("B2").Value = -124.3646

$ws.Range("A3").Value = 43.2397
$ws.Range("B3").Value = 43.2505

$ws.Range("A4").Value = -124.2333
$ws.Range("B4").Value = -124.2482

$ws.Range("A5").Value = 43.3455
$ws.Range("B5").Value = 43.3347
